# Updated cryptos list: new Price (D) and Volume(1h) (E) values per row.
# D-column values are written via a text-forced NumberFormat so that numeric-
# looking strings (e.g. "172.75") are preserved verbatim as text instead of
# being coerced into floating point numbers (matches original inlineStr cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '42.488.65'; E = '  +1.00%  ' },
    @{ Row = 3; D = '2.211.36'; E = '  -1.23%  ' },
    @{ Row = 4; D = $null; E = '  +0.01%  ' },
    @{ Row = 5; D = '240.73'; E = '  -0.92%  ' },
    @{ Row = 6; D = '0.615'; E = '  -1.43%  ' },
    @{ Row = 7; D = '74.76'; E = '  +0.65%  ' },
    @{ Row = 8; D = $null; E = '  -0.09%  ' },
    @{ Row = 9; D = '0.599'; E = '  -0.04%  ' },
    @{ Row = 10; D = '41.33'; E = '  -1.89%  ' },
    @{ Row = 11; D = '0.0923'; E = '  -3.13%  ' },
    @{ Row = 12; D = '54.92'; E = '  -2.84%  ' },
    @{ Row = 13; D = '6.86'; E = '  -0.88%  ' },
    @{ Row = 14; D = $null; E = '  -2.65%  ' },
    @{ Row = 15; D = '2.540.98'; E = '  -1.30%  ' },
    @{ Row = 16; D = '14.67'; E = '  +2.51%  ' },
    @{ Row = 17; D = '2.211.11'; E = '  -1.18%  ' },
    @{ Row = 18; D = '0.799'; E = '  -4.66%  ' },
    @{ Row = 19; D = '42.363.55'; E = '  +0.89%  ' },
    @{ Row = 20; D = '0.0000105'; E = '  -0.24%  ' },
    @{ Row = 21; D = '70.73'; E = '  -2.61%  ' },
    @{ Row = 22; D = '5.90'; E = '  -5.04%  ' },
    @{ Row = 23; D = '10.10'; E = '  -10.11%  ' },
    @{ Row = 24; D = '228.52'; E = '  -0.61%  ' },
    @{ Row = 25; D = $null; E = '  +2.61%  ' },
    @{ Row = 26; D = $null; E = '  +0.06%  ' },
    @{ Row = 27; D = '10.93'; E = '  -4.05%  ' },
    @{ Row = 28; D = '3.36'; E = '  -7.33%  ' },
    @{ Row = 29; D = $null; E = '  -2.26%  ' },
    @{ Row = 30; D = $null; E = '  -0.91%  ' },
    @{ Row = 31; D = '172.75'; E = $null },
    @{ Row = 32; D = '20.15'; E = $null },
    @{ Row = 33; D = '33.61'; E = '  +11.82%  ' },
    @{ Row = 34; D = '0.0792'; E = '  -1.27%  ' },
    @{ Row = 35; D = '5.38'; E = '  -4.08%  ' },
    @{ Row = 36; D = $null; E = '  -2.38%  ' },
    @{ Row = 37; D = $null; E = '  +3.06%  ' },
    @{ Row = 38; D = '0.108'; E = '  -2.87%  ' },
    @{ Row = 39; D = '0.0320'; E = '  +5.33%  ' },
    @{ Row = 40; D = '12.58'; E = '  -4.18%  ' },
    @{ Row = 41; D = '2.12'; E = '  -0.29%  ' },
    @{ Row = 42; D = '5.48'; E = '  -3.78%  ' },
    @{ Row = 43; D = '60.60'; E = '  -6.39%  ' },
    @{ Row = 44; D = '0.195'; E = '  -1.78%  ' },
    @{ Row = 45; D = '8.55'; E = '  -1.83%  ' },
    @{ Row = 46; D = $null; E = '  -2.28%  ' },
    @{ Row = 47; D = '99.39'; E = '  -4.23%  ' },
    @{ Row = 48; D = $null; E = '  -2.95%  ' },
    @{ Row = 49; D = '2.29'; E = '  -1.75%  ' },
    @{ Row = 50; D = '1.14'; E = '  -2.83%  ' },
    @{ Row = 51; D = '0.423'; E = '  +14.23%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D" + $u.Row)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.ClearFormats()
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
